$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("A1").Value = "season_ending_year"
$ws.Range("B1").Value = "lg"
$ws.Range("C1").Value = "blk_per_game"

# Copy the header style from A1 (bold/border/centered) onto the new C1 cell
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Update data rows with new values
$ws.Range("A2").Value = 1971
$ws.Range("B2").Value = "NBA"

$ws.Range("A3").Value = 1975
$ws.Range("B3").Value = "ABA"

$ws.Range("A4").Value = 1974
$ws.Range("B4").Value = "ABA"

$ws.Range("A5").Value = 1972
$ws.Range("B5").Value = "BAA"

# Clear column C for data rows (no data beyond header in column C)
$ws.Range("C2:C5").ClearContents()
